$d = $word.ActiveDocument

# Locate the paragraph that ends the "visibilidad" bullet list item; the
# three new bullet points about "reputacion" get appended right after it.
$rng = $d.Content
$found = $rng.Find.Execute(
    "No se podrá cambiar el tipo de visibilidad de la publicación.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

$anchor = $rng.Paragraphs.First

$newTexts = @(
    "Creamos la tabla Reputacion con 5 tipos de reputaciones, la reputación de un usuario se calcula como la suma de todas las estrellas recibidas dividido la cantidad de calificaciones.",
    "Todos los usuarios nuevos no poseen reputación, cada vez que alguien califica, se modifica la reputación del usuario calificado.",
    "Al final de la migración calculamos todas las reputaciones de los usuarios migrados, basándonos en las calificaciones migradas."
)

$prevPara = $anchor
foreach ($t in $newTexts) {
    # Inserting a new paragraph after an existing list item copies its
    # pStyle ("Prrafodelista") and numPr (list numId/ilvl) automatically.
    $prevPara.Range.InsertParagraphAfter()
    $newIndex = $prevPara.Index + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $newPara.Range.Text = $t
    $prevPara = $newPara
}
